$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add a new "2022-Q4" sheet, cloned from the current "2022-Q3"
#    sheet (same fund / same layout / same styles), placed right
#    before it, then overwrite the quarter-specific figures.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

# Columns D:G on this sheet are stored as text, so force a text
# number format before writing the numeric-looking strings.
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "113.53"
$q4.Range("E2").Value = "92.20"
$q4.Range("F2").Value = "3.63"
$q4.Range("G2").Value = "4.1211"
$q4.Range("H2").Value = 10

# ------------------------------------------------------------------
# 2. Insert a fresh top row in the "总计" (total) sheet for the new
#    quarter, pushing the older rows down by one.
# ------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows("2:2").Insert()

# The insert operation leaves stray formatting on the new row (e.g. it
# picks up bold from the header); strip it so B2:D2 stay plain like
# every other data row, then re-apply the bold+bordered "index column"
# look to just the A2 cell by copying formatting from the (shifted) A3.
$tot.Range("A2:D2").ClearFormats()
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q4"
$tot.Range("C2").Value = 1
$tot.Range("D2").Value = 4.12

# Renumber the running index in column A for every row that shifted
# down (it is a simple 0-based row counter).
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4

# ------------------------------------------------------------------
# 3. Restore the originally-active tab ("2021-Q2", the last sheet)
#    since inserting/copying sheets above moved the selection.
# ------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
